$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

# The text-like columns (Date, Time, Weekday, Week) must be stored as text,
# not auto-converted by Excel into dates/numbers. Temporarily force a text
# number format, assign the values, then clear the formatting again so the
# new row ends up with the same (default/no) style as the other data rows.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-12"
$ws.Cells.Item($row, 2).Value = "10:30:45"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "02"

$textRange.ClearFormats()

$ws.Cells.Item($row, 5).Value = 127154
$ws.Cells.Item($row, 6).Value = 143647
$ws.Cells.Item($row, 7).Value = 169185
$ws.Cells.Item($row, 8).Value = 159713
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142819
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193061
$ws.Cells.Item($row, 14).Value = 115508
$ws.Cells.Item($row, 15).Value = 45847
$ws.Cells.Item($row, 16).Value = 28496
$ws.Cells.Item($row, 17).Value = 65166
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48771
$ws.Cells.Item($row, 20).Value = -1
